$d = $word.ActiveDocument
$t = $d.Tables.Item(1)

$newValues = @(
"69-20=","26+35=","9+30=","93-3=","54-38=",
"59+9=","82-52=","91-48=","26+14=","16+31=",
"19+9=","95-58=","31-14=","42+44=","65+13=",
"99-57=","14-9=","28+6=","71-14=","97-79=",
"30+30=","55+41=","39+12=","65+29=","80-65=",
"55-36=","7+67=","0+38=","60-6=","46+39=",
"91+2=","61-12=","52+8=","18+48=","31+37=",
"64+25=","92-53=","64-9=","34+2=","98-61=",
"20+10=","63-26=","67-35=","11-4=","95-51=",
"36-12=","11+38=","84-16=","52-46=","73-19=",
"51-28=","75-11=","20+8=","75-2=","42-9=",
"19+34=","14+41=","50-26=","26+33=","65+34=",
"91-2=","28+68=","98-7=","21+48=","56+23=",
"32+42=","81-53=","23+22=","34+16=","16+56=",
"3+88=","39-5=","9+9=","37+2=","25+55=",
"25+30=","56+26=","78-17=","83-12=","22+32=",
"69+23=","23-7=","60-30=","18+74=","92-51=",
"72-48=","25+28=","91-10=","67+13=","86-33=",
"12+49=","2+10=","8+44=","44+24=","10-4=",
"46+53=","40+42=","33-13=","16+55=","37-3="
)

$rows = 20
$cols = 5
$idx = 0
for ($r = 1; $r -le $rows; $r++) {
    for ($c = 1; $c -le $cols; $c++) {
        $cell = $t.Cell($r, $c)
        $cell.Range.Text = $newValues[$idx]
        $idx = $idx + 1
    }
}
